$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 40
$ws1.Range("F4").Value = 1527
$ws1.Range("F5").Value = 226
$ws1.Range("F6").Value = 47
$ws1.Range("F7").Value = 591
$ws1.Range("F8").Value = 10003
$ws1.Range("F9").Value = 171
$ws1.Range("F11").Value = 243
$ws1.Range("F13").Value = 379
$ws1.Range("F14").Value = 6926
$ws1.Range("F15").Value = 1090
$ws1.Range("F16").Value = 645
$ws1.Range("F17").Value = 54
$ws1.Range("F18").Value = 206

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 40
$ws4.Range("F4").Value = 1527
$ws4.Range("F5").Value = 226
$ws4.Range("F7").Value = 47
$ws4.Range("F8").Value = 592
$ws4.Range("F11").Value = 10003
$ws4.Range("F12").Value = 171
$ws4.Range("F14").Value = 243
$ws4.Range("F16").Value = 379
$ws4.Range("F17").Value = 6926
$ws4.Range("F18").Value = 1090
$ws4.Range("F19").Value = 645
$ws4.Range("F20").Value = 54
$ws4.Range("F21").Value = 206

$wb.Save()
